$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers to Excel's parser need to be
# pre-formatted as Text so they keep their original string representation
# (leading/trailing zeros, no scientific notation) instead of being coerced to a number.
$textCells = @("D5", "D6", "D7", "D9", "D12", "D13", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D28", "D31", "D34", "D35", "D36", "D37", "D42", "D44", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.266.18"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.116.57"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "580.06"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "173.67"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "6.53"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "0.0000249"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "36.91"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D15").Value = "3.630.14"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "67.211.82"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "7.11"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "3.113.27"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "16.66"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "492.41"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("B21").Value = "Polygon"
$ws.Range("C21").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D21").Value = "0.706"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "7.91"
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "13.22"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").Value = "10.60"
$ws.Range("E26").Value = "  +5.67%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "7.98"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "28.41"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "0.0₃0943"
$ws.Range("E33").Value = "  -6.19%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "5.88"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "47.28"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D42").Value = "388.95"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "2.810.84"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").Value = "  -7.58%  "
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D48").Value = "25.13"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").Value = "6.73"
$ws.Range("E51").Value = "  -1.99%  "
